$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.990.10"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "2.757.21"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.63"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.95"
$ws.Range("E6").Value = "  -1.51%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  -3.51%  "

$ws.Range("E9").Value = "  -3.91%  "

$ws.Range("E10").Value = "  +3.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -15.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  -3.33%  "

$ws.Range("D13").Value = "3.245.94"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.91"
$ws.Range("E14").Value = "  -3.39%  "

$ws.Range("D15").Value = "63.611.16"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("E16").Value = "  -5.82%  "

$ws.Range("D17").Value = "2.763.07"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("E19").Value = "  -3.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "356.09"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  -5.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.527"
$ws.Range("E23").Value = "  -4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.89"
$ws.Range("E24").Value = "  -4.13%  "

$ws.Range("E25").Value = "  -2.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.56"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").Value = "0.0₃0901"
$ws.Range("E28").Value = "  -7.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -5.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.43"
$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.10"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -3.65%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("E37").Value = "  -3.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -2.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "350.02"
$ws.Range("E39").Value = "  +2.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.25"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("E41").Value = "  -2.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.09"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.43"
$ws.Range("E43").Value = "  -5.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.76"
$ws.Range("E44").Value = "  -4.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0589"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.93"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.633"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.19%  "
